$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) values for rows 2-49
# NumberFormat "@" (Text) is applied before writing numeric-looking Price
# strings so Excel keeps them as text instead of auto-converting to numbers.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "31.204.68"
$ws.Cells.Item(2, 5).Value = "  +2.15%  "
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.997.39"
$ws.Cells.Item(3, 5).Value = "  +6.17%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.000"
$ws.Cells.Item(4, 5).Value = "  +0.09%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "0.7798"
$ws.Cells.Item(5, 5).Value = "  +65.08%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "256.04"
$ws.Cells.Item(6, 5).Value = "  +3.93%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.9995"
$ws.Cells.Item(7, 5).Value = "  -0.02%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3479"
$ws.Cells.Item(8, 5).Value = "  +20.50%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "28.24"
$ws.Cells.Item(9, 5).Value = "  +28.08%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "44.30"
$ws.Cells.Item(10, 5).Value = "  +2.95%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.06999"
$ws.Cells.Item(11, 5).Value = "  +7.06%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.8518"
$ws.Cells.Item(12, 5).Value = "  +12.39%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.08184"
$ws.Cells.Item(13, 5).Value = "  +4.63%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "1.996.91"
$ws.Cells.Item(14, 5).Value = "  +6.19%  "
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "100.58"
$ws.Cells.Item(15, 5).Value = "  -0.50%  "
$ws.Cells.Item(16, 5).Value = "  +6.71%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "15.44"
$ws.Cells.Item(17, 5).Value = "  +17.19%  "
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "274.16"
$ws.Cells.Item(18, 5).Value = "  -3.82%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "31.220.09"
$ws.Cells.Item(19, 5).Value = "  +2.26%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "5.889"
$ws.Cells.Item(20, 5).Value = "  +9.69%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.000007910"
$ws.Cells.Item(21, 5).Value = "  +5.07%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "2.259.49"
$ws.Cells.Item(22, 5).Value = "  +6.79%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "1.000"
$ws.Cells.Item(23, 5).Value = "  +0.11%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "1.000"
$ws.Cells.Item(24, 5).Value = "  +0.13%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "7.107"
$ws.Cells.Item(25, 5).Value = "  +11.16%  "
$ws.Cells.Item(26, 5).Value = "  +10.03%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "164.42"
$ws.Cells.Item(27, 5).Value = "  +0.81%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.1479"
$ws.Cells.Item(28, 5).Value = "  +52.47%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "19.91"
$ws.Cells.Item(29, 5).Value = "  +4.40%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.321"
$ws.Cells.Item(30, 5).Value = "  +21.33%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.609"
$ws.Cells.Item(31, 5).Value = "  +7.68%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.609"
$ws.Cells.Item(32, 5).Value = "  +8.09%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.358"
$ws.Cells.Item(33, 5).Value = "  +2.46%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "4.400"
$ws.Cells.Item(34, 5).Value = "  +4.99%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.05232"
$ws.Cells.Item(35, 5).Value = "  +8.12%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.236"
$ws.Cells.Item(36, 5).Value = "  +9.34%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.7761"
$ws.Cells.Item(37, 5).Value = "  +11.64%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.757"
$ws.Cells.Item(38, 5).Value = "  -0.55%  "
$ws.Cells.Item(39, 5).Value = "  +4.55%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "2.902"
$ws.Cells.Item(40, 5).Value = "  +1.30%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "6.663"
$ws.Cells.Item(41, 5).Value = "  +4.80%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "79.38"
$ws.Cells.Item(42, 5).Value = "  +3.79%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "2.142"
$ws.Cells.Item(43, 5).Value = "  +8.21%  "
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.4674"
$ws.Cells.Item(44, 5).Value = "  +9.84%  "
$ws.Cells.Item(45, 5).Value = "  +4.29%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.8505"
$ws.Cells.Item(46, 5).Value = "  +2.55%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.9995"
$ws.Cells.Item(47, 5).Value = "  +0.02%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "7.708"
$ws.Cells.Item(48, 5).Value = "  +9.24%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "9.908"
$ws.Cells.Item(49, 5).Value = "  +0.92%  "

# Rows 50-51: Decentraland now ranks above Elrond (order swapped)
$ws.Cells.Item(50, 2).Value = "Decentraland"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.4323"
$ws.Cells.Item(50, 5).Value = "  +9.36%  "

$ws.Cells.Item(51, 2).Value = "Elrond"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "36.83"
$ws.Cells.Item(51, 5).Value = "  +4.85%  "
